$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6 for the "Trial: Attenuated" entry (shifts old rows 6-20 down to 7-21)
$ws.Rows.Item(6).Insert()

# Append new rows 22-24 at the bottom (Catch trial / white noise / no white noise)
$ws.Range("A22").Value = 555
$ws.Range("B22").Value = "Catch trial"

$ws.Range("A23").Value = 444
$ws.Range("B23").Value = "white noise"

$ws.Range("A24").Value = 400
$ws.Range("B24").Value = "no white noise"

# Fill in the newly inserted row 6: A=120, B="Trial: Attenuated"
$ws.Range("A6").Value = 120
$ws.Range("B6").Value = "Trial: Attenuated"

# Append new rows 25-26 at the bottom (attenuated correct rejection / false alarm)
$ws.Range("A25").Value = 22
$ws.Range("B25").Value = "correct rejection attenuated"

$ws.Range("A26").Value = 23
$ws.Range("B26").Value = "False alarm attenuated"

# Update the two relabeled rows (originally "correct rejection" / "False alarm",
# now shifted to rows 20 and 21) to the "(no object)" variants
$ws.Range("B20").Value = "correct rejection (no object)"
$ws.Range("B21").Value = "False alarm (no object)"

# Update the selected cell to match the final state of the workbook
$ws.Range("B29").Select()
